# Update the "想去人数" (want-to-go count) column F values on the
# "展览" and "全部类型" sheets to reflect newly generated scrape data.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    2  = 8467
    3  = 8154
    4  = 140
    9  = 145
    11 = 245
    12 = 735
    14 = 4228
    20 = 119
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
